# Imported Tuning as option in settings
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Typography sheet: add a new "Widget Wildcard Characters" column (J)
# ---------------------------------------------------------------------
$wsTypography = $wb.Worksheets.Item("Typography")
$wsTypography.Range("J3").Value = "Widget Wildcard Characters"
$wsTypography.Range("J4:J16").Value = ""

# ---------------------------------------------------------------------
# 2. Translation sheet: drop the obsolete "SingleUseId220" warning text
#    (row 219) and append the two new Tuning-related rows at the end.
# ---------------------------------------------------------------------
$wsTranslation = $wb.Worksheets.Item("Translation")

# Remove the row - everything below shifts up by one.
$wsTranslation.Rows.Item(219).Delete()

# New row 245: "Start Tuning Cycle" button text.
$wsTranslation.Range("B245").Value = "SingleUseId251"
$wsTranslation.Range("C245").Value = "Default"
$wsTranslation.Range("D245").Value = "Center"
$wsTranslation.Range("E245").Value = "LTR"
$wsTranslation.Range("F245").Value = "Start Tuning Cycle"
$wsTranslation.Range("G245").Value = "New Text"
$wsTranslation.Range("I245").Value = "New Text"
$wsTranslation.Range("J245").Value = "New Text"
$wsTranslation.Range("K245").Value = "New Text"
$wsTranslation.Range("L245").Value = "New Text"
$wsTranslation.Range("M245").Value = "New Text"
$wsTranslation.Range("N245").Value = "New Text"

# New row 246: "Tuning" label, fully translated.
$wsTranslation.Range("B246").Value = "Tuning"
$wsTranslation.Range("C246").Value = "Large"
$wsTranslation.Range("D246").Value = "Center"
$wsTranslation.Range("E246").Value = "LTR"
$wsTranslation.Range("F246").Value = "Tuning"
$wsTranslation.Range("G246").Value = "调音"
$wsTranslation.Range("H246").Value = "LargeZHS"
$wsTranslation.Range("I246").Value = "Réglage"
$wsTranslation.Range("J246").Value = "Sintonización"
$wsTranslation.Range("K246").Value = "Tuning"
$wsTranslation.Range("L246").Value = "Messa a punto"
$wsTranslation.Range("M246").Value = "Tuning"
$wsTranslation.Range("N246").Value = "настройка"

Write-Host "Tuning option imported."
